$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "10uF Ceramic Capacitor, 25V"
$ws.Range("D3").Value = "C96446"

$ws.Range("D4").Value = "C1591"

$ws.Range("D5").Value = "C2762594"

$ws.Range("C6").Value = "SMD,7.2x6.6mm"

$ws.Range("A9").Value = "56.2k Resistor, 0.1%"
$ws.Range("D9").Value = "C705784"

$ws.Range("A10").Value = "10k Resistor, 0.1%"
$ws.Range("D10").Value = "C95204"

$ws.Columns.Item(2).ColumnWidth = 85
